$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 20 de Marzo de 2020 a las 08:46"
$ws.Cells.Item(32, 1).Value2 = "Pakistan"
$ws.Cells.Item(32, 2).Value2 = 454
$ws.Cells.Item(32, 3).Value2 = 0
$ws.Cells.Item(32, 4).Value2 = 13
$ws.Cells.Item(32, 5).Value2 = 438
$ws.Cells.Item(32, 6).Value2 = 0
$ws.Cells.Item(32, 7).Value2 = 1
$ws.Cells.Item(32, 8).Value2 = 3

$ws.Cells.Item(34, 1).Value2 = "Polonia"
$ws.Cells.Item(34, 2).Value2 = 365
$ws.Cells.Item(34, 3).Value2 = 10
$ws.Cells.Item(34, 4).Value2 = 13
$ws.Cells.Item(34, 5).Value2 = 347
$ws.Cells.Item(34, 6).Value2 = 3
$ws.Cells.Item(34, 7).Value2 = 0
$ws.Cells.Item(34, 8).Value2 = 5

$ws.Cells.Item(35, 1).Value2 = "Turquia"
$ws.Cells.Item(35, 2).Value2 = 359
$ws.Cells.Item(35, 3).Value2 = 0
$ws.Cells.Item(35, 4).Value2 = 0
$ws.Cells.Item(35, 5).Value2 = 355
$ws.Cells.Item(35, 6).Value2 = 0
$ws.Cells.Item(35, 7).Value2 = 0
$ws.Cells.Item(35, 8).Value2 = 4

$ws.Cells.Item(62, 1).Value2 = "Armenia"
$ws.Cells.Item(62, 2).Value2 = 136
$ws.Cells.Item(62, 3).Value2 = 14
$ws.Cells.Item(62, 4).Value2 = 1
$ws.Cells.Item(62, 5).Value2 = 135
$ws.Cells.Item(62, 6).Value2 = 2
$ws.Cells.Item(62, 7).Value2 = 0
$ws.Cells.Item(62, 8).Value2 = 0

$ws.Cells.Item(63, 1).Value2 = "Colombia"
$ws.Cells.Item(63, 2).Value2 = 128
$ws.Cells.Item(63, 3).Value2 = 20
$ws.Cells.Item(63, 4).Value2 = 1
$ws.Cells.Item(63, 5).Value2 = 127
$ws.Cells.Item(63, 6).Value2 = 0
$ws.Cells.Item(63, 7).Value2 = 0
$ws.Cells.Item(63, 8).Value2 = 0

$ws.Cells.Item(64, 1).Value2 = "Argentina"
$ws.Cells.Item(64, 2).Value2 = 128
$ws.Cells.Item(64, 3).Value2 = 0
$ws.Cells.Item(64, 4).Value2 = 3
$ws.Cells.Item(64, 5).Value2 = 122
$ws.Cells.Item(64, 6).Value2 = 0
$ws.Cells.Item(64, 7).Value2 = 0
$ws.Cells.Item(64, 8).Value2 = 3

$ws.Cells.Item(65, 1).Value2 = "Eslovaquia"
$ws.Cells.Item(65, 2).Value2 = 124
$ws.Cells.Item(65, 3).Value2 = 0
$ws.Cells.Item(65, 4).Value2 = 0
$ws.Cells.Item(65, 5).Value2 = 124
$ws.Cells.Item(65, 6).Value2 = 2
$ws.Cells.Item(65, 7).Value2 = 0
$ws.Cells.Item(65, 8).Value2 = 0

$ws.Cells.Item(66, 1).Value2 = "Serbia"
$ws.Cells.Item(66, 2).Value2 = 118
$ws.Cells.Item(66, 3).Value2 = 15
$ws.Cells.Item(66, 4).Value2 = 2
$ws.Cells.Item(66, 5).Value2 = 116
$ws.Cells.Item(66, 6).Value2 = 4
$ws.Cells.Item(66, 7).Value2 = 0
$ws.Cells.Item(66, 8).Value2 = 0

$ws.Cells.Item(86, 1).Value2 = "Bielorrusia"
$ws.Cells.Item(86, 2).Value2 = 51
$ws.Cells.Item(86, 3).Value2 = 0
$ws.Cells.Item(86, 4).Value2 = 15
$ws.Cells.Item(86, 5).Value2 = 36
$ws.Cells.Item(86, 6).Value2 = 0
$ws.Cells.Item(86, 7).Value2 = 0
$ws.Cells.Item(86, 8).Value2 = 0

$ws.Cells.Item(117, 1).Value2 = "Montenegro"
$ws.Cells.Item(117, 2).Value2 = 13
$ws.Cells.Item(117, 3).Value2 = 0
$ws.Cells.Item(117, 4).Value2 = 0
$ws.Cells.Item(117, 5).Value2 = 13
$ws.Cells.Item(117, 6).Value2 = 0
$ws.Cells.Item(117, 7).Value2 = 0
$ws.Cells.Item(117, 8).Value2 = 0

$ws.Cells.Item(118, 1).Value2 = "Paraguay"
$ws.Cells.Item(118, 2).Value2 = 13
$ws.Cells.Item(118, 3).Value2 = 0
$ws.Cells.Item(118, 4).Value2 = 0
$ws.Cells.Item(118, 5).Value2 = 13
$ws.Cells.Item(118, 6).Value2 = 1
$ws.Cells.Item(118, 7).Value2 = 0
$ws.Cells.Item(118, 8).Value2 = 0

$ws.Cells.Item(119, 1).Value2 = "Maldivas"
$ws.Cells.Item(119, 2).Value2 = 13
$ws.Cells.Item(119, 3).Value2 = 0
$ws.Cells.Item(119, 4).Value2 = 0
$ws.Cells.Item(119, 5).Value2 = 13
$ws.Cells.Item(119, 6).Value2 = 0
$ws.Cells.Item(119, 7).Value2 = 0
$ws.Cells.Item(119, 8).Value2 = 0

$ws.Cells.Item(120, 1).Value2 = "Camerun"
$ws.Cells.Item(120, 2).Value2 = 13
$ws.Cells.Item(120, 3).Value2 = 0
$ws.Cells.Item(120, 4).Value2 = 2
$ws.Cells.Item(120, 5).Value2 = 11
$ws.Cells.Item(120, 6).Value2 = 0
$ws.Cells.Item(120, 7).Value2 = 0
$ws.Cells.Item(120, 8).Value2 = 0

$ws.Cells.Item(130, 1).Value2 = "Togo"
$ws.Cells.Item(130, 2).Value2 = 9
$ws.Cells.Item(130, 3).Value2 = 8
$ws.Cells.Item(130, 4).Value2 = 0
$ws.Cells.Item(130, 5).Value2 = 9
$ws.Cells.Item(130, 6).Value2 = 0
$ws.Cells.Item(130, 7).Value2 = 0
$ws.Cells.Item(130, 8).Value2 = 0

$ws.Cells.Item(132, 1).Value2 = "Costa de Marfil"
$ws.Cells.Item(132, 2).Value2 = 9
$ws.Cells.Item(132, 3).Value2 = 0
$ws.Cells.Item(132, 4).Value2 = 1
$ws.Cells.Item(132, 5).Value2 = 8
$ws.Cells.Item(132, 6).Value2 = 0
$ws.Cells.Item(132, 7).Value2 = 0
$ws.Cells.Item(132, 8).Value2 = 0

$ws.Cells.Item(134, 1).Value2 = "Mauricio"
$ws.Cells.Item(134, 2).Value2 = 7
$ws.Cells.Item(134, 3).Value2 = 0
$ws.Cells.Item(134, 4).Value2 = 0
$ws.Cells.Item(134, 5).Value2 = 7
$ws.Cells.Item(134, 6).Value2 = 0
$ws.Cells.Item(134, 7).Value2 = 0
$ws.Cells.Item(134, 8).Value2 = 0

$ws.Cells.Item(136, 1).Value2 = "Puerto Rico"
$ws.Cells.Item(136, 2).Value2 = 6
$ws.Cells.Item(136, 3).Value2 = 0
$ws.Cells.Item(136, 4).Value2 = 0
$ws.Cells.Item(136, 5).Value2 = 6
$ws.Cells.Item(136, 6).Value2 = 0
$ws.Cells.Item(136, 7).Value2 = 0
$ws.Cells.Item(136, 8).Value2 = 0

$ws.Cells.Item(137, 1).Value2 = "Guinea Ecuatorial"
$ws.Cells.Item(137, 2).Value2 = 6
$ws.Cells.Item(137, 3).Value2 = 0
$ws.Cells.Item(137, 4).Value2 = 0
$ws.Cells.Item(137, 5).Value2 = 6
$ws.Cells.Item(137, 6).Value2 = 0
$ws.Cells.Item(137, 7).Value2 = 0
$ws.Cells.Item(137, 8).Value2 = 0

$ws.Cells.Item(138, 1).Value2 = "Mongolia"
$ws.Cells.Item(138, 2).Value2 = 6
$ws.Cells.Item(138, 3).Value2 = 0
$ws.Cells.Item(138, 4).Value2 = 0
$ws.Cells.Item(138, 5).Value2 = 6
$ws.Cells.Item(138, 6).Value2 = 0
$ws.Cells.Item(138, 7).Value2 = 0
$ws.Cells.Item(138, 8).Value2 = 0

$ws.Cells.Item(139, 1).Value2 = "Tanzania"
$ws.Cells.Item(139, 2).Value2 = 6
$ws.Cells.Item(139, 3).Value2 = 0
$ws.Cells.Item(139, 4).Value2 = 0
$ws.Cells.Item(139, 5).Value2 = 6
$ws.Cells.Item(139, 6).Value2 = 0
$ws.Cells.Item(139, 7).Value2 = 0
$ws.Cells.Item(139, 8).Value2 = 0

$ws.Cells.Item(140, 1).Value2 = "Seychelles"
$ws.Cells.Item(140, 2).Value2 = 6
$ws.Cells.Item(140, 3).Value2 = 0
$ws.Cells.Item(140, 4).Value2 = 0
$ws.Cells.Item(140, 5).Value2 = 6
$ws.Cells.Item(140, 6).Value2 = 0
$ws.Cells.Item(140, 7).Value2 = 0
$ws.Cells.Item(140, 8).Value2 = 0

$ws.Cells.Item(141, 1).Value2 = "Barbados"
$ws.Cells.Item(141, 2).Value2 = 5
$ws.Cells.Item(141, 3).Value2 = 0
$ws.Cells.Item(141, 4).Value2 = 0
$ws.Cells.Item(141, 5).Value2 = 5
$ws.Cells.Item(141, 6).Value2 = 0
$ws.Cells.Item(141, 7).Value2 = 0
$ws.Cells.Item(141, 8).Value2 = 0

$ws.Cells.Item(143, 1).Value2 = "Aruba"
$ws.Cells.Item(143, 2).Value2 = 5
$ws.Cells.Item(143, 3).Value2 = 0
$ws.Cells.Item(143, 4).Value2 = 1
$ws.Cells.Item(143, 5).Value2 = 4
$ws.Cells.Item(143, 6).Value2 = 0
$ws.Cells.Item(143, 7).Value2 = 0
$ws.Cells.Item(143, 8).Value2 = 0

$ws.Cells.Item(144, 1).Value2 = "Mayotte"
$ws.Cells.Item(144, 2).Value2 = 4
$ws.Cells.Item(144, 3).Value2 = 0
$ws.Cells.Item(144, 4).Value2 = 0
$ws.Cells.Item(144, 5).Value2 = 4
$ws.Cells.Item(144, 6).Value2 = 0
$ws.Cells.Item(144, 7).Value2 = 0
$ws.Cells.Item(144, 8).Value2 = 0

$ws.Cells.Item(145, 1).Value2 = "Congo"
$ws.Cells.Item(145, 2).Value2 = 3
$ws.Cells.Item(145, 3).Value2 = 0
$ws.Cells.Item(145, 4).Value2 = 0
$ws.Cells.Item(145, 5).Value2 = 3
$ws.Cells.Item(145, 6).Value2 = 0
$ws.Cells.Item(145, 7).Value2 = 0
$ws.Cells.Item(145, 8).Value2 = 0

$ws.Cells.Item(146, 1).Value2 = "San Martin (Parte Francesa)"
$ws.Cells.Item(146, 2).Value2 = 3
$ws.Cells.Item(146, 3).Value2 = 0
$ws.Cells.Item(146, 4).Value2 = 0
$ws.Cells.Item(146, 5).Value2 = 3
$ws.Cells.Item(146, 6).Value2 = 0
$ws.Cells.Item(146, 7).Value2 = 0
$ws.Cells.Item(146, 8).Value2 = 0

$ws.Cells.Item(147, 1).Value2 = "San Bartolome"
$ws.Cells.Item(147, 2).Value2 = 3
$ws.Cells.Item(147, 3).Value2 = 0
$ws.Cells.Item(147, 4).Value2 = 0
$ws.Cells.Item(147, 5).Value2 = 3
$ws.Cells.Item(147, 6).Value2 = 0
$ws.Cells.Item(147, 7).Value2 = 0
$ws.Cells.Item(147, 8).Value2 = 0

$ws.Cells.Item(148, 1).Value2 = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(148, 2).Value2 = 3
$ws.Cells.Item(148, 3).Value2 = 0
$ws.Cells.Item(148, 4).Value2 = 0
$ws.Cells.Item(148, 5).Value2 = 3
$ws.Cells.Item(148, 6).Value2 = 0
$ws.Cells.Item(148, 7).Value2 = 0
$ws.Cells.Item(148, 8).Value2 = 0

$ws.Cells.Item(149, 1).Value2 = "Bahamas"
$ws.Cells.Item(149, 2).Value2 = 3
$ws.Cells.Item(149, 3).Value2 = 0
$ws.Cells.Item(149, 4).Value2 = 0
$ws.Cells.Item(149, 5).Value2 = 3
$ws.Cells.Item(149, 6).Value2 = 0
$ws.Cells.Item(149, 7).Value2 = 0
$ws.Cells.Item(149, 8).Value2 = 0

$ws.Cells.Item(151, 1).Value2 = "Gabon"
$ws.Cells.Item(151, 2).Value2 = 3
$ws.Cells.Item(151, 3).Value2 = 0
$ws.Cells.Item(151, 4).Value2 = 0
$ws.Cells.Item(151, 5).Value2 = 3
$ws.Cells.Item(151, 6).Value2 = 0
$ws.Cells.Item(151, 7).Value2 = 0
$ws.Cells.Item(151, 8).Value2 = 0

$ws.Cells.Item(153, 1).Value2 = "Islas Caimanes"
$ws.Cells.Item(153, 2).Value2 = 3
$ws.Cells.Item(153, 3).Value2 = 0
$ws.Cells.Item(153, 4).Value2 = 0
$ws.Cells.Item(153, 5).Value2 = 2
$ws.Cells.Item(153, 6).Value2 = 0
$ws.Cells.Item(153, 7).Value2 = 0
$ws.Cells.Item(153, 8).Value2 = 1

$ws.Cells.Item(154, 1).Value2 = "Bermudas"
$ws.Cells.Item(154, 2).Value2 = 2
$ws.Cells.Item(154, 3).Value2 = 0
$ws.Cells.Item(154, 4).Value2 = 0
$ws.Cells.Item(154, 5).Value2 = 2
$ws.Cells.Item(154, 6).Value2 = 0
$ws.Cells.Item(154, 7).Value2 = 0
$ws.Cells.Item(154, 8).Value2 = 0

$ws.Cells.Item(155, 1).Value2 = "Benin"
$ws.Cells.Item(155, 2).Value2 = 2
$ws.Cells.Item(155, 3).Value2 = 0
$ws.Cells.Item(155, 4).Value2 = 0
$ws.Cells.Item(155, 5).Value2 = 2
$ws.Cells.Item(155, 6).Value2 = 0
$ws.Cells.Item(155, 7).Value2 = 0
$ws.Cells.Item(155, 8).Value2 = 0

$ws.Cells.Item(156, 1).Value2 = "Liberia"
$ws.Cells.Item(156, 2).Value2 = 2
$ws.Cells.Item(156, 3).Value2 = 0
$ws.Cells.Item(156, 4).Value2 = 0
$ws.Cells.Item(156, 5).Value2 = 2
$ws.Cells.Item(156, 6).Value2 = 0
$ws.Cells.Item(156, 7).Value2 = 0
$ws.Cells.Item(156, 8).Value2 = 0

$ws.Cells.Item(157, 1).Value2 = "Mauritania"
$ws.Cells.Item(157, 2).Value2 = 2
$ws.Cells.Item(157, 3).Value2 = 0
$ws.Cells.Item(157, 4).Value2 = 0
$ws.Cells.Item(157, 5).Value2 = 2
$ws.Cells.Item(157, 6).Value2 = 0
$ws.Cells.Item(157, 7).Value2 = 0
$ws.Cells.Item(157, 8).Value2 = 0

$ws.Cells.Item(158, 1).Value2 = "Republica de Africa Central"
$ws.Cells.Item(158, 2).Value2 = 2
$ws.Cells.Item(158, 3).Value2 = 1
$ws.Cells.Item(158, 4).Value2 = 0
$ws.Cells.Item(158, 5).Value2 = 2
$ws.Cells.Item(158, 6).Value2 = 0
$ws.Cells.Item(158, 7).Value2 = 0
$ws.Cells.Item(158, 8).Value2 = 0

$ws.Cells.Item(159, 1).Value2 = "Butan"
$ws.Cells.Item(159, 2).Value2 = 2
$ws.Cells.Item(159, 3).Value2 = 1
$ws.Cells.Item(159, 4).Value2 = 0
$ws.Cells.Item(159, 5).Value2 = 2
$ws.Cells.Item(159, 6).Value2 = 0
$ws.Cells.Item(159, 7).Value2 = 0
$ws.Cells.Item(159, 8).Value2 = 0

$ws.Cells.Item(160, 1).Value2 = "Zambia"
$ws.Cells.Item(160, 2).Value2 = 2
$ws.Cells.Item(160, 3).Value2 = 0
$ws.Cells.Item(160, 4).Value2 = 0
$ws.Cells.Item(160, 5).Value2 = 2
$ws.Cells.Item(160, 6).Value2 = 0
$ws.Cells.Item(160, 7).Value2 = 0
$ws.Cells.Item(160, 8).Value2 = 0

$ws.Cells.Item(161, 1).Value2 = "Groenlandia"
$ws.Cells.Item(161, 2).Value2 = 2
$ws.Cells.Item(161, 3).Value2 = 0
$ws.Cells.Item(161, 4).Value2 = 0
$ws.Cells.Item(161, 5).Value2 = 2
$ws.Cells.Item(161, 6).Value2 = 0
$ws.Cells.Item(161, 7).Value2 = 0
$ws.Cells.Item(161, 8).Value2 = 0

$ws.Cells.Item(162, 1).Value2 = "Santa Lucia"
$ws.Cells.Item(162, 2).Value2 = 2
$ws.Cells.Item(162, 3).Value2 = 0
$ws.Cells.Item(162, 4).Value2 = 0
$ws.Cells.Item(162, 5).Value2 = 2
$ws.Cells.Item(162, 6).Value2 = 0
$ws.Cells.Item(162, 7).Value2 = 0
$ws.Cells.Item(162, 8).Value2 = 0

$ws.Cells.Item(163, 1).Value2 = "Nueva Caledonia"
$ws.Cells.Item(163, 2).Value2 = 2
$ws.Cells.Item(163, 3).Value2 = 0
$ws.Cells.Item(163, 4).Value2 = 0
$ws.Cells.Item(163, 5).Value2 = 2
$ws.Cells.Item(163, 6).Value2 = 0
$ws.Cells.Item(163, 7).Value2 = 0
$ws.Cells.Item(163, 8).Value2 = 0

$ws.Cells.Item(164, 1).Value2 = "Haiti"
$ws.Cells.Item(164, 2).Value2 = 2
$ws.Cells.Item(164, 3).Value2 = 2
$ws.Cells.Item(164, 4).Value2 = 0
$ws.Cells.Item(164, 5).Value2 = 2
$ws.Cells.Item(164, 6).Value2 = 0
$ws.Cells.Item(164, 7).Value2 = 0
$ws.Cells.Item(164, 8).Value2 = 0

$ws.Cells.Item(165, 1).Value2 = "Sudan"
$ws.Cells.Item(165, 2).Value2 = 2
$ws.Cells.Item(165, 3).Value2 = 0
$ws.Cells.Item(165, 4).Value2 = 0
$ws.Cells.Item(165, 5).Value2 = 1
$ws.Cells.Item(165, 6).Value2 = 0
$ws.Cells.Item(165, 7).Value2 = 0
$ws.Cells.Item(165, 8).Value2 = 1

$ws.Cells.Item(166, 1).Value2 = "Isla de Man"
$ws.Cells.Item(166, 2).Value2 = 1
$ws.Cells.Item(166, 3).Value2 = 0
$ws.Cells.Item(166, 4).Value2 = 0
$ws.Cells.Item(166, 5).Value2 = 1
$ws.Cells.Item(166, 6).Value2 = 0
$ws.Cells.Item(166, 7).Value2 = 0
$ws.Cells.Item(166, 8).Value2 = 0

$ws.Cells.Item(167, 1).Value2 = "San Vicente y las Granadinas"
$ws.Cells.Item(167, 2).Value2 = 1
$ws.Cells.Item(167, 3).Value2 = 0
$ws.Cells.Item(167, 4).Value2 = 0
$ws.Cells.Item(167, 5).Value2 = 1
$ws.Cells.Item(167, 6).Value2 = 0
$ws.Cells.Item(167, 7).Value2 = 0
$ws.Cells.Item(167, 8).Value2 = 0

$ws.Cells.Item(168, 1).Value2 = "Republica de Yibuti"
$ws.Cells.Item(168, 2).Value2 = 1
$ws.Cells.Item(168, 3).Value2 = 0
$ws.Cells.Item(168, 4).Value2 = 0
$ws.Cells.Item(168, 5).Value2 = 1
$ws.Cells.Item(168, 6).Value2 = 0
$ws.Cells.Item(168, 7).Value2 = 0
$ws.Cells.Item(168, 8).Value2 = 0

$ws.Cells.Item(169, 1).Value2 = "Surinam"
$ws.Cells.Item(169, 2).Value2 = 1
$ws.Cells.Item(169, 3).Value2 = 0
$ws.Cells.Item(169, 4).Value2 = 0
$ws.Cells.Item(169, 5).Value2 = 1
$ws.Cells.Item(169, 6).Value2 = 0
$ws.Cells.Item(169, 7).Value2 = 0
$ws.Cells.Item(169, 8).Value2 = 0

$ws.Cells.Item(170, 1).Value2 = "Guinea"
$ws.Cells.Item(170, 2).Value2 = 1
$ws.Cells.Item(170, 3).Value2 = 0
$ws.Cells.Item(170, 4).Value2 = 0
$ws.Cells.Item(170, 5).Value2 = 1
$ws.Cells.Item(170, 6).Value2 = 0
$ws.Cells.Item(170, 7).Value2 = 0
$ws.Cells.Item(170, 8).Value2 = 0

$ws.Cells.Item(172, 1).Value2 = "San Martin (Parte Holandesa)"
$ws.Cells.Item(172, 2).Value2 = 1
$ws.Cells.Item(172, 3).Value2 = 0
$ws.Cells.Item(172, 4).Value2 = 0
$ws.Cells.Item(172, 5).Value2 = 1
$ws.Cells.Item(172, 6).Value2 = 0
$ws.Cells.Item(172, 7).Value2 = 0
$ws.Cells.Item(172, 8).Value2 = 0

$ws.Cells.Item(173, 1).Value2 = "Angola"
$ws.Cells.Item(173, 2).Value2 = 1
$ws.Cells.Item(173, 3).Value2 = 1
$ws.Cells.Item(173, 4).Value2 = 0
$ws.Cells.Item(173, 5).Value2 = 1
$ws.Cells.Item(173, 6).Value2 = 0
$ws.Cells.Item(173, 7).Value2 = 0
$ws.Cells.Item(173, 8).Value2 = 0

$ws.Cells.Item(174, 1).Value2 = "Montserrat"
$ws.Cells.Item(174, 2).Value2 = 1
$ws.Cells.Item(174, 3).Value2 = 0
$ws.Cells.Item(174, 4).Value2 = 0
$ws.Cells.Item(174, 5).Value2 = 1
$ws.Cells.Item(174, 6).Value2 = 0
$ws.Cells.Item(174, 7).Value2 = 0
$ws.Cells.Item(174, 8).Value2 = 0

$ws.Cells.Item(175, 1).Value2 = "Gambia"
$ws.Cells.Item(175, 2).Value2 = 1
$ws.Cells.Item(175, 3).Value2 = 0
$ws.Cells.Item(175, 4).Value2 = 0
$ws.Cells.Item(175, 5).Value2 = 1
$ws.Cells.Item(175, 6).Value2 = 0
$ws.Cells.Item(175, 7).Value2 = 0
$ws.Cells.Item(175, 8).Value2 = 0

$ws.Cells.Item(176, 1).Value2 = "El Salvador"
$ws.Cells.Item(176, 2).Value2 = 1
$ws.Cells.Item(176, 3).Value2 = 0
$ws.Cells.Item(176, 4).Value2 = 0
$ws.Cells.Item(176, 5).Value2 = 1
$ws.Cells.Item(176, 6).Value2 = 0
$ws.Cells.Item(176, 7).Value2 = 0
$ws.Cells.Item(176, 8).Value2 = 0

$ws.Cells.Item(177, 1).Value2 = "Nicaragua"
$ws.Cells.Item(177, 2).Value2 = 1
$ws.Cells.Item(177, 3).Value2 = 0
$ws.Cells.Item(177, 4).Value2 = 0
$ws.Cells.Item(177, 5).Value2 = 1
$ws.Cells.Item(177, 6).Value2 = 0
$ws.Cells.Item(177, 7).Value2 = 0
$ws.Cells.Item(177, 8).Value2 = 0

$ws.Cells.Item(178, 1).Value2 = "Suazilandia"
$ws.Cells.Item(178, 2).Value2 = 1
$ws.Cells.Item(178, 3).Value2 = 0
$ws.Cells.Item(178, 4).Value2 = 0
$ws.Cells.Item(178, 5).Value2 = 1
$ws.Cells.Item(178, 6).Value2 = 0
$ws.Cells.Item(178, 7).Value2 = 0
$ws.Cells.Item(178, 8).Value2 = 0

$ws.Cells.Item(179, 1).Value2 = "Republica del Chad"
$ws.Cells.Item(179, 2).Value2 = 1
$ws.Cells.Item(179, 3).Value2 = 0
$ws.Cells.Item(179, 4).Value2 = 0
$ws.Cells.Item(179, 5).Value2 = 1
$ws.Cells.Item(179, 6).Value2 = 0
$ws.Cells.Item(179, 7).Value2 = 0
$ws.Cells.Item(179, 8).Value2 = 0

$ws.Cells.Item(180, 1).Value2 = "Antigua y Barbuda"
$ws.Cells.Item(180, 2).Value2 = 1
$ws.Cells.Item(180, 3).Value2 = 0
$ws.Cells.Item(180, 4).Value2 = 0
$ws.Cells.Item(180, 5).Value2 = 1
$ws.Cells.Item(180, 6).Value2 = 0
$ws.Cells.Item(180, 7).Value2 = 0
$ws.Cells.Item(180, 8).Value2 = 0

$ws.Cells.Item(181, 1).Value2 = "Niger"
$ws.Cells.Item(181, 2).Value2 = 1
$ws.Cells.Item(181, 3).Value2 = 0
$ws.Cells.Item(181, 4).Value2 = 0
$ws.Cells.Item(181, 5).Value2 = 1
$ws.Cells.Item(181, 6).Value2 = 0
$ws.Cells.Item(181, 7).Value2 = 0
$ws.Cells.Item(181, 8).Value2 = 0

$ws.Cells.Item(182, 1).Value2 = "Fiyi"
$ws.Cells.Item(182, 2).Value2 = 1
$ws.Cells.Item(182, 3).Value2 = 0
$ws.Cells.Item(182, 4).Value2 = 0
$ws.Cells.Item(182, 5).Value2 = 1
$ws.Cells.Item(182, 6).Value2 = 0
$ws.Cells.Item(182, 7).Value2 = 0
$ws.Cells.Item(182, 8).Value2 = 0

$ws.Cells.Item(183, 1).Value2 = "Santa Sede"
$ws.Cells.Item(183, 2).Value2 = 1
$ws.Cells.Item(183, 3).Value2 = 0
$ws.Cells.Item(183, 4).Value2 = 0
$ws.Cells.Item(183, 5).Value2 = 1
$ws.Cells.Item(183, 6).Value2 = 0
$ws.Cells.Item(183, 7).Value2 = 0
$ws.Cells.Item(183, 8).Value2 = 0

$ws.Cells.Item(184, 1).Value2 = "Somalia"
$ws.Cells.Item(184, 2).Value2 = 1
$ws.Cells.Item(184, 3).Value2 = 0
$ws.Cells.Item(184, 4).Value2 = 0
$ws.Cells.Item(184, 5).Value2 = 1
$ws.Cells.Item(184, 6).Value2 = 0
$ws.Cells.Item(184, 7).Value2 = 0
$ws.Cells.Item(184, 8).Value2 = 0

